$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2 = 3
    3 = 6
    4 = 4
    5 = 10
    6 = 8
    7 = 10
    8 = 4
    9 = 7
    10 = 3
    11 = 6
    12 = 7
    13 = 5
    14 = 6
    15 = 3
    16 = 7
    17 = 4
    18 = 5
    19 = 4
    20 = 5
    21 = 2
    22 = 9
    23 = 5
    24 = 3
    25 = 1
    26 = 3
    27 = 7
    28 = 0
    29 = 1
    30 = 4
    31 = 5
    32 = 4
    33 = 6
    34 = 5
    35 = 2
    36 = 1
    37 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
